$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = '30.392.71'
$ws.Cells.Item(2, 5).Value = '  +2.32%  '
$ws.Cells.Item(3, 4).Value = '2.108.84'
$ws.Cells.Item(3, 5).Value = '  +0.51%  '
$ws.Cells.Item(4, 5).Value = '  -0.18%  '
Set-TextValue $ws.Cells.Item(5, 4) '344.95'
$ws.Cells.Item(5, 5).Value = '  +0.63%  '
Set-TextValue $ws.Cells.Item(6, 4) '1.005'
$ws.Cells.Item(6, 5).Value = '  -0.14%  '
Set-TextValue $ws.Cells.Item(7, 4) '0.5235'
$ws.Cells.Item(7, 5).Value = '  +2.05%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.4447'
$ws.Cells.Item(8, 5).Value = '  +1.07%  '
$ws.Cells.Item(9, 5).Value = '  +3.19%  '
Set-TextValue $ws.Cells.Item(10, 4) '0.09401'
$ws.Cells.Item(10, 5).Value = '  +2.49%  '
$ws.Cells.Item(11, 5).Value = '  +0.35%  '
Set-TextValue $ws.Cells.Item(12, 4) '24.98'
$ws.Cells.Item(12, 5).Value = '  +0.59%  '
Set-TextValue $ws.Cells.Item(13, 4) '8.679'
$ws.Cells.Item(13, 5).Value = '  +6.17%  '
Set-TextValue $ws.Cells.Item(14, 4) '6.947'
$ws.Cells.Item(14, 5).Value = '  +2.99%  '
$ws.Cells.Item(15, 4).Value = '2.088.34'
$ws.Cells.Item(15, 5).Value = '  -0.64%  '
Set-TextValue $ws.Cells.Item(16, 4) '101.98'
$ws.Cells.Item(16, 5).Value = '  +2.43%  '
Set-TextValue $ws.Cells.Item(17, 4) '0.00001164'
$ws.Cells.Item(17, 5).Value = '  +1.21%  '
Set-TextValue $ws.Cells.Item(19, 4) '21.25'
$ws.Cells.Item(19, 5).Value = '  +0.91%  '
Set-TextValue $ws.Cells.Item(20, 4) '0.06722'
$ws.Cells.Item(20, 5).Value = '  +1.21%  '
Set-TextValue $ws.Cells.Item(21, 4) '6.360'
$ws.Cells.Item(21, 5).Value = '  +3.02%  '
Set-TextValue $ws.Cells.Item(22, 4) '1.005'
$ws.Cells.Item(22, 5).Value = '  -0.21%  '
$ws.Cells.Item(23, 4).Value = '30.430.42'
$ws.Cells.Item(23, 5).Value = '  +2.33%  '
Set-TextValue $ws.Cells.Item(24, 4) '12.65'
$ws.Cells.Item(24, 5).Value = '  +0.53%  '
Set-TextValue $ws.Cells.Item(25, 4) '2.297'
$ws.Cells.Item(25, 5).Value = '  -0.44%  '
Set-TextValue $ws.Cells.Item(26, 4) '22.04'
$ws.Cells.Item(26, 5).Value = '  +0.99%  '
Set-TextValue $ws.Cells.Item(27, 4) '162.93'
$ws.Cells.Item(27, 5).Value = '  +0.17%  '
Set-TextValue $ws.Cells.Item(28, 4) '2.529'
$ws.Cells.Item(28, 5).Value = '  +0.67%  '
Set-TextValue $ws.Cells.Item(29, 4) '133.93'
$ws.Cells.Item(29, 5).Value = '  +1.13%  '
Set-TextValue $ws.Cells.Item(30, 4) '1.151'
$ws.Cells.Item(30, 5).Value = '  +1.99%  '
Set-TextValue $ws.Cells.Item(31, 4) '1.753'
$ws.Cells.Item(31, 5).Value = '  +7.46%  '
$ws.Cells.Item(32, 5).Value = '  +0.91%  '
Set-TextValue $ws.Cells.Item(33, 4) '6.819'
$ws.Cells.Item(33, 5).Value = '  +13.26%  '
Set-TextValue $ws.Cells.Item(34, 4) '6.274'
$ws.Cells.Item(34, 5).Value = '  +1.96%  '
Set-TextValue $ws.Cells.Item(35, 4) '3.923'
$ws.Cells.Item(35, 5).Value = '  -1.02%  '
Set-TextValue $ws.Cells.Item(36, 4) '10.36'
$ws.Cells.Item(36, 5).Value = '  +1.35%  '
Set-TextValue $ws.Cells.Item(37, 4) '0.02637'
$ws.Cells.Item(37, 5).Value = '  +2.89%  '
Set-TextValue $ws.Cells.Item(38, 4) '0.06801'
$ws.Cells.Item(38, 5).Value = '  +1.64%  '
Set-TextValue $ws.Cells.Item(39, 4) '0.7065'
$ws.Cells.Item(39, 5).Value = '  +3.25%  '
$ws.Cells.Item(40, 5).Value = '  +5.37%  '
Set-TextValue $ws.Cells.Item(41, 4) '12.61'
$ws.Cells.Item(41, 5).Value = '  +1.84%  '
$ws.Cells.Item(42, 5).Value = '  +0.21%  '
$ws.Cells.Item(43, 5).Value = '  +2.86%  '
Set-TextValue $ws.Cells.Item(44, 4) '14.42'
$ws.Cells.Item(44, 5).Value = '  +1.56%  '
Set-TextValue $ws.Cells.Item(45, 4) '2.365'
$ws.Cells.Item(45, 5).Value = '  +3.03%  '
Set-TextValue $ws.Cells.Item(46, 4) '1.005'
$ws.Cells.Item(46, 5).Value = '  -0.09%  '
Set-TextValue $ws.Cells.Item(47, 4) '1.391'
$ws.Cells.Item(47, 5).Value = '  +19.57%  '
Set-TextValue $ws.Cells.Item(48, 4) '3.648'
$ws.Cells.Item(48, 5).Value = '  +1.05%  '
Set-TextValue $ws.Cells.Item(49, 4) '0.00000000355'
$ws.Cells.Item(49, 5).Value = '  +6.20%  '
Set-TextValue $ws.Cells.Item(50, 4) '1.214'
$ws.Cells.Item(50, 5).Value = '  +9.85%  '
Set-TextValue $ws.Cells.Item(51, 4) '1.221'
$ws.Cells.Item(51, 5).Value = '  +0.34%  '
